$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value/type corrections (NaN marker <-> numeric) across existing rows ---
$ws.Range("I9").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("CS16").Value = "NaN"
$ws.Range("CS17").Value = "NaN"
$ws.Range("I18").Value = "NaN"
$ws.Range("L19").Value = "NaN"
$ws.Range("DA19").Value = "NaN"
$ws.Range("CW20").Value = "NaN"
$ws.Range("CW21").Value = "NaN"
$ws.Range("CW22").Value = "NaN"
$ws.Range("BG23").Value = "NaN"
$ws.Range("CW23").Value = "NaN"
$ws.Range("CW26").Value = 1
$ws.Range("CW27").Value = 1
$ws.Range("CM31").Value = 1
$ws.Range("AK33").Value = "NaN"
$ws.Range("CV38").Value = "NaN"
$ws.Range("DN74").Value = "NaN"
$ws.Range("DN75").Value = "NaN"
$ws.Range("CM80").Value = 17
$ws.Range("CM81").Value = 20
$ws.Range("CM82").Value = "NaN"
$ws.Range("DQ89").Value = "NaN"
$ws.Range("H90").Value = 3
$ws.Range("CF93").Value = "NaN"
$ws.Range("CF107").Value = "NaN"

# --- Append new data row 174 ---
$ws.Range("A174").Value = 44068
$ws.Range("B174").Value = 562128
$ws.Range("C174").Value = 2691
$ws.Range("D174").Value = 72495
$ws.Range("E174").Value = 63255
$ws.Range("F174").Value = 195137
$ws.Range("G174").Value = 24158
$ws.Range("H174").Value = 3136
$ws.Range("I174").Value = 2474
$ws.Range("J174").Value = 5179
$ws.Range("K174").Value = 4408
$ws.Range("L174").Value = 8385
$ws.Range("M174").Value = 3675
$ws.Range("N174").Value = 18187
$ws.Range("O174").Value = 20336
$ws.Range("P174").Value = 4608
$ws.Range("Q174").Value = 3644
$ws.Range("R174").Value = 11624
$ws.Range("S174").Value = 6562
$ws.Range("T174").Value = 13177
$ws.Range("U174").Value = 9548
$ws.Range("V174").Value = 2599
$ws.Range("W174").Value = 968
$ws.Range("X174").Value = 5002
$ws.Range("Y174").Value = 14703
$ws.Range("Z174").Value = 10858
$ws.Range("AA174").Value = 6016
$ws.Range("AB174").Value = 44394
$ws.Range("AC174").Value = 914
$ws.Range("AD174").Value = 172
$ws.Range("AE174").Value = 229
$ws.Range("AF174").Value = 441
$ws.Range("AG174").Value = 74
$ws.Range("AH174").Value = 36
$ws.Range("AI174").Value = 237
$ws.Range("AJ174").Value = 1941
$ws.Range("AK174").Value = 2487
$ws.Range("AL174").Value = 35654
$ws.Range("AM174").Value = 6147
$ws.Range("AN174").Value = 2390
$ws.Range("AO174").Value = 34903
$ws.Range("AP174").Value = 852
$ws.Range("AQ174").Value = 19695
$ws.Range("AR174").Value = 1418
$ws.Range("AS174").Value = 6670
$ws.Range("AT174").Value = 1423
$ws.Range("AU174").Value = 1544
$ws.Range("AV174").Value = 3541
$ws.Range("AW174").Value = 1483
$ws.Range("AX174").Value = 925
$ws.Range("AY174").Value = 2459
$ws.Range("AZ174").Value = 2574
$ws.Range("BA174").Value = 41630
$ws.Range("BB174").Value = 11403
$ws.Range("BC174").Value = 2069
$ws.Range("BD174").Value = 7025
$ws.Range("BE174").Value = 3233
$ws.Range("BF174").Value = 274
$ws.Range("BG174").Value = 1375
$ws.Range("BH174").Value = 2558
$ws.Range("BI174").Value = 728
$ws.Range("BJ174").Value = 1959
$ws.Range("BK174").Value = 7857
$ws.Range("BL174").Value = 7820
$ws.Range("BM174").Value = 7571
$ws.Range("BN174").Value = 13639
$ws.Range("BO174").Value = 1863
$ws.Range("BP174").Value = 788
$ws.Range("BQ174").Value = 5799
$ws.Range("BR174").Value = 5484
$ws.Range("BS174").Value = 6127
$ws.Range("BT174").Value = 1302
$ws.Range("BU174").Value = 1394
$ws.Range("BV174").Value = 2377
$ws.Range("BW174").Value = 2719
$ws.Range("BX174").Value = 729
$ws.Range("BY174").Value = 3935
$ws.Range("BZ174").Value = 2253
$ws.Range("CA174").Value = 1136
$ws.Range("CB174").Value = 640
$ws.Range("CC174").Value = 1866
$ws.Range("CD174").Value = 1797
$ws.Range("CE174").Value = 1025
$ws.Range("CF174").Value = 838
$ws.Range("CG174").Value = 4325
$ws.Range("CH174").Value = 1150
$ws.Range("CI174").Value = 1106
$ws.Range("CJ174").Value = 1125
$ws.Range("CK174").Value = 1423
$ws.Range("CL174").Value = 1326
$ws.Range("CM174").Value = 1375
$ws.Range("CN174").Value = 1064
$ws.Range("CO174").Value = 1013
$ws.Range("CP174").Value = 1055
$ws.Range("CQ174").Value = 542
$ws.Range("CR174").Value = 2869
$ws.Range("CS174").Value = 901
$ws.Range("CT174").Value = 772
$ws.Range("CU174").Value = 703
$ws.Range("CV174").Value = 1212
$ws.Range("CW174").Value = 1053
$ws.Range("CX174").Value = 577
$ws.Range("CY174").Value = 696
$ws.Range("CZ174").Value = 797
$ws.Range("DA174").Value = 1078
$ws.Range("DB174").Value = 886
$ws.Range("DC174").Value = 1018
$ws.Range("DD174").Value = 787
$ws.Range("DE174").Value = 311
$ws.Range("DF174").Value = 325
$ws.Range("DG174").Value = 646
$ws.Range("DH174").Value = 548
$ws.Range("DI174").Value = 390
$ws.Range("DJ174").Value = 529
$ws.Range("DK174").Value = 316
$ws.Range("DL174").Value = 562
$ws.Range("DM174").Value = 693
$ws.Range("DN174").Value = 503
$ws.Range("DO174").Value = 471
$ws.Range("DP174").Value = 354
$ws.Range("DQ174").Value = 508
$ws.Range("DR174").Value = 112030
$ws.Range("DS174").Value = 236552
$ws.Range("DT174").Value = 8717
$ws.Range("DU174").Value = 101530
$ws.Range("DV174").Value = 65654
$ws.Range("DW174").Value = 24300
$ws.Range("DX174").Value = 7852

# --- Update active cell selection to the new last row ---
$ws.Range("A174").Select()
